$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.600824356079102
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.645496129989624
$ws.Range("D1").Value = 1.507542490959167
$ws.Range("E1").Value = 1.164753437042236
